$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Step 1: restructure cumcontrol1/cumcontrol2/dickpic/boosters sheets ---
# Delete the old "dickpic" sheet (its content is fully superseded by the old cumcontrol2 content).
$wb.Worksheets.Item("dickpic").Delete()

# The old "cumcontrol2" sheet content becomes the new "dickpic" sheet (relabeled rows + new text).
$wsDick = $wb.Worksheets.Item("cumcontrol2")
$wsDick.Range("A2").Value = 'dpppv2'
$wsDick.Range("B2").Value = 'okay you just made me do something, give me a sec'
$wsDick.Range("C2").Value = 'LEVERAGE variant.'
$wsDick.Range("A3").Value = 'dpppv1'
$wsDick.Range("B3").Value = 'you can''t just send me that and expect me to do nothing about it, hold on...'
$wsDick.Range("C3").Value = 'LEVERAGE. WAIT 1-2 min then send PPV.'
$wsDick.Range("A4").Value = 'dprapport2'
$wsDick.Range("B4").Value = 'woah I wasn''t expecting that but... damn okay'
$wsDick.Range("C4").Value = 'DURING RAPPORT variant.'
$wsDick.Range("A5").Value = 'dprapport1'
$wsDick.Range("B5").Value = 'omg you don''t waste time huh, that''s actually really hot ngl'
$wsDick.Range("C5").Value = 'DURING RAPPORT.'
$wsDick.Range("A6").Value = 'dpsext2'
$wsDick.Range("B6").Value = 'oh fuck that is... damn, I need to show you something rn'
$wsDick.Range("C6").Value = 'DURING SEXTING variant.'
$wsDick.Range("A7").Value = 'dpsext1'
$wsDick.Range("B7").Value = 'fuck okay that''s... wow, you have no idea what that just did to me'
$wsDick.Range("C7").Value = 'DURING SEXTING.'
$wsDick.Name = "dickpic"

# The old "cumcontrol1" sheet is edited in place and renamed to "cumcontrol".
$wsCum = $wb.Worksheets.Item("cumcontrol1")
$wsCum.Range("B2").Value = 'don''t you dare finish before you see this'
$wsCum.Range("B3").Value = 'hold it, I want you to wait until you see what I''m about to send, trust me it''s worth the wait'
$wsCum.Range("C3").Value = 'DELAY. Send PPV.'
$wsCum.Range("B4").Value = 'wait for me, I want us to finish together, open this first'
$wsCum.Range("C4").Value = 'SYNC variant.'
$wsCum.Range("B5").Value = 'I''m so close too, cum with me... but you need to see this first'
$wsCum.Range("B6").Value = 'hold it, not yet... I need you to last longer for me'
$wsCum.Range("B7").Value = 'don''t cum yet, I''m not done with you'
$wsCum.Range("C7").Value = 'CONTROL.'
$wsCum.Name = "cumcontrol"

# The "boosters" sheet keeps its name; only text content changes.
$wsBoost = $wb.Worksheets.Item("boosters")
$wsBoost.Range("B2").Value = 'I literally can''t focus on anything else rn'
$wsBoost.Range("C2").Value = 'BOOSTER.'
$wsBoost.Range("B3").Value = 'more'
$wsBoost.Range("B4").Value = 'my hands are shaking'
$wsBoost.Range("B5").Value = 'I can''t think straight rn'
$wsBoost.Range("B6").Value = 'you have no idea what you''re doing to me'
$wsBoost.Range("B7").Value = 'don''t stop'
$wsBoost.Range("B8").Value = 'I''m so wet rn because of you'
$wsBoost.Range("B9").Value = 'fuckkk'

# --- Step 2: text-only edits on the remaining conversation-script sheets ---

$ws = $wb.Worksheets.Item("free1")
$ws.Range("B2").Value = 'I''ll keep it then, but just know I made it thinking about you'
$ws.Range("B3").Value = 'I literally spent time making this because of what YOU said to me, this wasn''t random content, this was for you'
$ws.Range("B4").Value = 'free? I don''t just hand this out, you have to earn the best stuff from me'
$ws.Range("B5").Value = 'baby I already gave you one for free remember? this one is on another level and you know it'
$ws.Range("C5").Value = 'REMINDER. Still → Step 2.'

$ws = $wb.Worksheets.Item("free2")
$ws.Range("B2").Value = 'okay I''ll hold onto it, but it''s got your name on it'
$ws.Range("B3").Value = 'I recorded this because of our conversation, specifically because of YOU, it took effort and I did it for you'
$ws.Range("B4").Value = 'the best things aren''t free baby, especially not what I just did'
$ws.Range("B5").Value = 'you already got a free preview and this is ten times crazier, you know I don''t do free twice'

$ws = $wb.Worksheets.Item("nomoney1")
$ws.Range("B2").Value = 'it''s fine, you''re still here and that''s what matters to me right now'
$ws.Range("C2").Value = 'PROTECT.'
$ws.Range("B3").Value = 'send whatever you can, even something small, I just need you to see what you made me do'
$ws.Range("B4").Value = 'not even [small amount]? I really want you to be the one who sees this'
$ws.Range("C4").Value = 'TEST $3-5. Still no → Step 3.'
$ws.Range("B5").Value = 'I hear you baby, that''s okay'
$ws.Range("C5").Value = 'EMPATHY. Still → Step 2.'

$ws = $wb.Worksheets.Item("nomoney2")
$ws.Range("B2").Value = 'you being here talking to me is enough for now'
$ws.Range("B3").Value = 'just send whatever feels right, I can''t keep this from you'
$ws.Range("B4").Value = 'what about just [small amount]? I''d hate for you to miss this one'
$ws.Range("B5").Value = 'okay don''t stress about it'

$ws = $wb.Worksheets.Item("noppv1")
$ws.Range("B2").Value = 'send me whatever you want, literally anything, I just can''t keep this to myself, you need to see it'
$ws.Range("B3").Value = 'okay forget about money, this isn''t about that, I need you to see what you''re doing to me rn because I''ve never reacted like this'
$ws.Range("B4").Value = 'that''s fine, I''m not trying to sell you anything, I''m just enjoying this'

$ws = $wb.Worksheets.Item("noppv2")
$ws.Range("B2").Value = 'send anything, even the smallest amount, you need to see what you did to me'
$ws.Range("B3").Value = 'forget the money part for a sec, I want to share this moment with you because what you''re making me feel is real'
$ws.Range("B4").Value = 'no worries, I don''t care about that rn, this convo is what I care about'

$ws = $wb.Worksheets.Item("card1")
$ws.Range("B2").Value = 'figure it out soon because this mood I''m in right now isn''t gonna last and I want you to have it'
$ws.Range("B3").Value = 'maybe try a different card? I really don''t want you to miss this'
$ws.Range("C3").Value = 'ALTERNATIVE. Still → Step 3.'
$ws.Range("B4").Value = 'ugh that''s annoying, try again baby it usually works the second time'
$ws.Range("C4").Value = 'RETRY. Still → Step 2.'

$ws = $wb.Worksheets.Item("card2")
$ws.Range("B2").Value = 'baby fix it quick, I don''t keep stuff like this around forever'
$ws.Range("B3").Value = 'do you have another card? because you really don''t want to miss what I made for you'
$ws.Range("C3").Value = 'ALTERNATIVE. Still → Step 3.'
$ws.Range("B4").Value = 'that happens sometimes, try one more time for me'
$ws.Range("C4").Value = 'RETRY. Still → Step 2.'

$ws = $wb.Worksheets.Item("nosex1")
$ws.Range("B2").Value = 'alright I''ll behave, for now... no promises though'
$ws.Range("B3").Value = 'okay but I can''t help it, there''s something about you that keeps messing with my head'
$ws.Range("B4").Value = 'so tell me more about you... what do you do when you''re not making girls on the internet lose focus?'
$ws.Range("C4").Value = 'SUBTLE. → Step 3 later.'
$ws.Range("B5").Value = 'haha okay I got a little carried away, you''re too easy to talk to'
$ws.Range("C5").Value = 'RESPECT. Still → Step 2.'

$ws = $wb.Worksheets.Item("nosex2")
$ws.Range("B2").Value = 'fine I''ll stop, but don''t blame me when it happens again'
$ws.Range("B3").Value = 'I''m trying so hard to behave but you make it impossible'
$ws.Range("B4").Value = 'okay new topic, what''s the most adventurous thing you''ve ever done?'
$ws.Range("B5").Value = 'my bad I got ahead of myself, it''s your fault for being this fun to talk to'
$ws.Range("C5").Value = 'RESPECT. Still → Step 2.'

$ws = $wb.Worksheets.Item("offtopic1")
$ws.Range("B2").Value = 'okay wait I remember now, so like I was saying...'
$ws.Range("B3").Value = 'but wait you totally distracted me, I was about to say something and now you made me lose my train of thought'
$ws.Range("B4").Value = 'haha okay that''s actually funny'
$ws.Range("C4").Value = 'ACKNOWLEDGE. Adapt to what he said. → Step 2.'

$ws = $wb.Worksheets.Item("offtopic2")
$ws.Range("B2").Value = 'OKAY focus, where was I... oh yeah'
$ws.Range("B3").Value = 'wait no stop, you''re distracting me from what I was going to tell you'
$ws.Range("B4").Value = 'omg okay random but I love it'
$ws.Range("C4").Value = 'ACKNOWLEDGE. Adapt. → Step 2.'

$ws = $wb.Worksheets.Item("real1")
$ws.Range("B2").Value = 'I get why you''d think that, there''s a lot of fake stuff on here, but what we''ve been talking about felt real to me... didn''t it feel real to you?'
$ws.Range("B3").Value = 'ask me anything, literally anything about my life, I''m an open book, go ahead'
$ws.Range("B4").Value = 'lol do I sound like a robot to you? beep boop... send $5 for human verification haha I''m kidding'
$ws.Range("C4").Value = 'HUMOR. Still → Step 2.'

$ws = $wb.Worksheets.Item("real2")
$ws.Range("B2").Value = 'I know there''s a lot of bots out here but what we''ve been talking about... that felt different, you felt it too right?'
$ws.Range("B3").Value = 'test me then, ask me something only a real person would know, go ahead'
$ws.Range("B4").Value = 'wait you think I''m fake?? that might be the funniest thing anyone''s said to me today'

$ws = $wb.Worksheets.Item("voice1")
$ws.Range("A2").Value = 'Step3 Firm'
$ws.Range("B2").Value = 'I don''t do that on here but what I''m about to show you is better than any call, you''ll see'
$ws.Range("C2").Value = 'FIRM. No video calls.'
$ws.Range("A3").Value = 'Step2 Redirect'
$ws.Range("B3").Value = 'I have something way better for you though, trust me you''ll forget you even asked'
$ws.Range("C3").Value = 'REDIRECT. Still → Step 3.'
$ws.Range("A4").Value = 'Step1 Dodge'
$ws.Range("B4").Value = 'haha maybe one day if you earn it, but not yet'
$ws.Range("C4").Value = 'DODGE. Model does NOT do video calls. Still → Step 2.'

$ws = $wb.Worksheets.Item("voice2")
$ws.Range("A2").Value = 'Step3 Firm'
$ws.Range("B2").Value = 'I don''t do that here but trust me what I have is way better than my voice'
$ws.Range("C2").Value = 'FIRM.'
$ws.Range("A3").Value = 'Step2 Redirect'
$ws.Range("B3").Value = 'how about instead of a call I show you something that''ll actually blow your mind?'
$ws.Range("C3").Value = 'REDIRECT. Still → Step 3.'
$ws.Range("A4").Value = 'Step1 Dodge'
$ws.Range("B4").Value = 'hmmm you gotta earn that first'
$ws.Range("C4").Value = 'DODGE. Still → Step 2.'

$ws = $wb.Worksheets.Item("customyes1")
$ws.Range("B2").Value = 'trust me you won''t regret it, I made this one special'
$ws.Range("B3").Value = 'I have it and you''re gonna lose your mind... [price]'
$ws.Range("C3").Value = 'PRICE. Set based on content.'
$ws.Range("B4").Value = 'mmm you want that? I might have exactly what you''re thinking of'

$ws = $wb.Worksheets.Item("customyes2")
$ws.Range("B2").Value = 'you''re not gonna be able to stop watching this one'
$ws.Range("B3").Value = 'I made something just like that, [price] and it''s worth every penny'
$ws.Range("C3").Value = 'PRICE.'
$ws.Range("B4").Value = 'ohhh you have good taste, I think I know exactly what you need'

$ws = $wb.Worksheets.Item("customno1")
$ws.Range("B2").Value = 'trust me, I know what you need better than you do'
$ws.Range("B3").Value = 'what I have might be even crazier and literally no one else has seen it'
$ws.Range("B4").Value = 'I don''t have exactly that but I have something that''ll make you forget you even asked'
$ws.Range("C4").Value = 'REDIRECT. → Step 2.'

$ws = $wb.Worksheets.Item("customno2")
$ws.Range("B2").Value = 'just trust me on this one, you''ll thank me after'
$ws.Range("B3").Value = 'nobody has seen what I''m about to show you and I think it''s better than what you asked for'
$ws.Range("B4").Value = 'hmm not exactly that but what I DO have is gonna hit even harder'
$ws.Range("C4").Value = 'REDIRECT. → Step 2.'

$ws = $wb.Worksheets.Item("done1")
$ws.Range("B2").Value = 'next time you HAVE to wait for me, I have something insane planned for round 2'
$ws.Range("B3").Value = 'but I haven''t finished yet... you''re really gonna leave me like this?'
$ws.Range("C3").Value = 'RESCUE. Still no → Step 3.'
$ws.Range("B4").Value = 'fuck that''s so hot, you came because of me??'

$ws = $wb.Worksheets.Item("done2")
$ws.Range("B2").Value = 'okay but next time you hold it, because what I have planned is way crazier'
$ws.Range("B3").Value = 'wait but I''m not done yet, don''t you wanna watch me finish too?'
$ws.Range("B4").Value = 'already?? damn that''s hot'

$ws = $wb.Worksheets.Item("price1")
$ws.Range("B2").Value = 'alright, I''ll keep it then... but next time I''m in this mood you better be ready'
$ws.Range("B3").Value = 'fine, [lower price] but only because you''ve been making me feel some type of way, don''t tell anyone I did this'
$ws.Range("B4").Value = 'honestly? most guys couldn''t handle what I just recorded... I thought you were different'
$ws.Range("B5").Value = 'I''m literally dripping rn because of what you said and this mood isn''t gonna last, don''t miss it'
$ws.Range("B6").Value = 'babe that''s less than your morning coffee and I promise this is gonna keep you up way longer'
$ws.Range("C6").Value = 'REFRAME. Still no → Step 2.'

$ws = $wb.Worksheets.Item("price2")
$ws.Range("B2").Value = 'I''ll hold onto it for now, but I know you''re gonna come back for it'
$ws.Range("B3").Value = 'look, [lower price] because I actually want YOU to have this one, I''m not doing this for anyone else'
$ws.Range("B4").Value = 'maybe you''re just not ready for what I did... it''s a lot'
$ws.Range("B5").Value = 'this mood won''t last and I already picked YOU to share it with, don''t make me regret that'
$ws.Range("B6").Value = 'baby it''s literally nothing for what you''re about to see, trust me you''ll be thanking me after'

$ws = $wb.Worksheets.Item("discount1")
$ws.Range("B2").Value = 'okay I''ll just save this for someone who actually wants it then, their loss isn''t your gain though'
$ws.Range("C2").Value = 'TAKEAWAY.'
$ws.Range("B3").Value = 'ugh fine, [lower price] ONLY because I like you but don''t ever ask me this again'
$ws.Range("B4").Value = 'the guys who get my best stuff don''t ask for discounts, they know what they''re getting'
$ws.Range("B5").Value = 'haha negotiate? baby this isn''t a flea market... you already know this is worth every penny'

$ws = $wb.Worksheets.Item("discount2")
$ws.Range("B2").Value = 'no worries, I''ll keep it... someone else has been asking and they won''t hesitate'
$ws.Range("B3").Value = 'okay [lower price] and that''s ONLY because this convo has been different, first and last time'
$ws.Range("B4").Value = 'I don''t do this for just anyone, and the ones who get it never complain about the price after'
$ws.Range("B5").Value = 'a discount? babe do I look like I''m on sale? you know exactly what you''re getting'
